# Updates cryptos.xlsx price/volume columns (and a few re-ranked rows)
# to match the refreshed coinranking.com snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.676.41'
$ws.Range('E2').Value = '  -3.32%  '
$ws.Range('D3').Value = '2.607.77'
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '574.62'
$ws.Range('E5').Value = '  -4.02%  '
$ws.Range('D6').Value = '156.19'
$ws.Range('E6').Value = '  -1.94%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -3.35%  '
$ws.Range('D9').Value = '0.119'
$ws.Range('E9').Value = '  -6.62%  '
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('D11').Value = "'0.380"
$ws.Range('E11').Value = '  -5.34%  '
$ws.Range('E12').Value = '  -0.48%  '
$ws.Range('D13').Value = '28.22'
$ws.Range('E13').Value = '  -2.97%  '
$ws.Range('D14').Value = '3.076.85'
$ws.Range('E14').Value = '  -2.19%  '
$ws.Range('D15').Value = "'0.0000179"
$ws.Range('E15').Value = '  -8.39%  '
$ws.Range('D16').Value = '63.561.77'
$ws.Range('E16').Value = '  -3.41%  '
$ws.Range('D17').Value = '2.625.18'
$ws.Range('E17').Value = '  -0.90%  '
$ws.Range('D18').Value = '12.03'
$ws.Range('E18').Value = '  -4.73%  '
$ws.Range('D19').Value = '7.59'
$ws.Range('E19').Value = '  +1.30%  '
$ws.Range('D20').Value = '4.54'
$ws.Range('E20').Value = '  -5.75%  '
$ws.Range('D21').Value = '343.32'
$ws.Range('E21').Value = '  -2.71%  '
$ws.Range('E22').Value = '  +0.12%  '
$ws.Range('D23').Value = '67.56'
$ws.Range('E23').Value = '  -3.67%  '
$ws.Range('E24').Value = '  -0.66%  '
$ws.Range('D25').Value = "'0.0000109"
$ws.Range('E25').Value = '  -3.93%  '
$ws.Range('D26').Value = '9.18'
$ws.Range('E26').Value = '  -4.91%  '
$ws.Range('D27').Value = '585.04'
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('D28').Value = '1.58'
$ws.Range('E28').Value = '  -2.72%  '
$ws.Range('D29').Value = '0.162'
$ws.Range('E29').Value = '  -1.28%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('E31').Value = '  -3.10%  '
$ws.Range('E32').Value = '  -4.15%  '
$ws.Range('E33').Value = '  -4.19%  '
$ws.Range('D34').Value = '6.57'
$ws.Range('E34').Value = '  -2.43%  '
$ws.Range('D35').Value = '5.38'
$ws.Range('E35').Value = '  -2.88%  '
$ws.Range('D36').Value = '0.403'
$ws.Range('E36').Value = '  -4.74%  '
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').Value = '19.75'
$ws.Range('E37').Value = '  -4.17%  '
$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('D39').Value = '154.23'
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('E40').Value = '  -4.71%  '
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '2.54'
$ws.Range('E42').Value = '  +1.31%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').Value = '41.35'
$ws.Range('E43').Value = '  -3.46%  '
$ws.Range('D44').Value = '157.24'
$ws.Range('E44').Value = '  -2.89%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').Value = '23.81'
$ws.Range('E45').Value = '  +1.25%  '
$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '3.91'
$ws.Range('E46').Value = '  -4.87%  '
$ws.Range('D47').Value = '0.0589'
$ws.Range('E47').Value = '  -4.65%  '
$ws.Range('D48').Value = "'0.630"
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('E50').Value = '  -5.11%  '
$ws.Range('D51').Value = '18.81'
$ws.Range('E51').Value = '  -5.23%  '
